$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a single new row above row 2, pushing existing rows (old 2-18) down to (3-19)
$ws.Rows.Item(2).Insert()

# Row 3 now holds the old row-2 data ("On-screen debug text", 14) shifted down;
# overwrite it first with the new "Finish off debug rendering" task text, so the
# old "On-screen debug text" string becomes unused (and drops out of sharedStrings).
$ws.Range("A3").Value = "Engine"
$ws.Range("B3").Value = "Finish off debug rendering (text at least)"
$ws.Range("C3").Value = 4

# Fill the newly inserted (blank) row 2 with the "Engine refactoring" task.
$ws.Range("A2").Value = "Engine"
$ws.Range("B2").Value = "Engine refactoring (assume it is to be demoed"
$ws.Range("C2").Value = 28

# Select the new row as the active range (matches selection change in diff)
$ws.Range("A2:C2").Select()
